$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows (52-53) entered first so the new shared strings land at the
# same indices the authored workbook uses.
$ws.Range("A52").Value = "saveECModel.m"
$ws.Range("A53").Value = "enhanceGEM.m"
$ws.Range("C53").Value = "Probably not relevant anymore, remove?"

# Fill in the newly-added test-case notes in column B.
$ws.Range("B7").Value = "tc0001, tc0002"
$ws.Range("B6").Value = "tc0003, tc0004 - File download not covered by test cases but is tested manually and works."
$ws.Range("B8").Value = "tc0005"
$ws.Range("B3").Value = "tc0003, tc0004"
$ws.Range("B18").Value = "tc0006"
$ws.Range("B17").Value = "tc0007 - does not test download of the databases - this is tested in the manual workflows for Yeast-GEM and Human-GEM"
$ws.Range("B41").Value = "tc0008"
$ws.Range("B52").Value = "tc0009 - currently doesn't work"

# Row 53 reuses the existing "Not explicitly tested..." note.
$ws.Range("B53").Value = "Not explicitly tested - The code has existed for a long time"

# Match the saved selection / scroll position from the authored file.
$null = $ws.Range("A9").Select()
